# Review Log V 2.0
# Add a new review item (Review-17, "Git Hub Project Structure") on row 19,
# matching the layout/format already used for the Review-16 entry on row 18,
# and update the window's scroll/selection to reflect where the author was
# working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 19 values. A19 already holds "Review-17"; fill in the rest of that
# review row the same way row 18 ("Review-16") was filled in.
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = $ws.Range("C18").Text
$ws.Range("D19").Value = "Git Hub Project Structure"
$ws.Range("E19").Value = $ws.Range("E18").Text
$ws.Range("F19").Value = $ws.Range("F18").Text
$ws.Range("G19").Value = $ws.Range("G18").Text
$ws.Range("H19").Value = "The Project repository  needs to be organized (documents must have unified  format DocName_ID without any Version number (this one will be a tag on the tool) )"

# ---------------------------------------------------------------------------
# Formatting to match row 18's look: centered horizontally/vertically, no
# wrap for C:G, and the wrapped "What to Change" style (copied straight from
# H18) for H19.
# ---------------------------------------------------------------------------
$ws.Range("C19:G19").HorizontalAlignment = -4108
$ws.Range("C19:G19").VerticalAlignment = -4108
$ws.Range("C19:G19").WrapText = $false

$ws.Range("H18").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 19 needs to grow to fit the new "What to Change" text.
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------------
# Window state: scroll the grid and move the active selection to where the
# author left off after adding the new row.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F19").Select()
